# feat: add 2022-Q1 data
#
# The former "总计" (summary) sheet becomes the new "2022-Q1" sheet (reusing
# its sheetId/relationship slot), populated with the 2022-Q1 fund holdings.
# A brand-new "总计" sheet is appended at the end with the rebuilt summary
# table (now including the 2022-Q1 row at the top).

# Excel auto-converts numeric-looking strings ("001208", "14.51", ...) to
# real numbers on assignment, which both loses the leading zeros on fund
# codes and introduces binary floating-point noise on the decimal figures.
# Force those cells to Text first, then drop back to the default ("Normal")
# style so no stray number-format style is left behind on the cell.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")

# Grab fully-styled cells (bold/centered/bordered header, bold/centered row
# index) from an existing per-quarter sheet so the new columns/cells can
# copy their exact formatting.
$headerStyleSource = $wb.Worksheets.Item("2021-Q4").Range("B1")
$indexStyleSource = $wb.Worksheets.Item("2021-Q4").Range("A2")

# --- Step 1: the old "总计" sheet is reborn as "2022-Q1" -------------------
$summary.Name = "2022-Q1"
$q1 = $summary

# Clear out the old 4-row "总计" table before writing the new layout.
$q1.Cells.Clear()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerStyleSource.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Row index column (A) style, copied from the same source sheet.
$indexStyleSource.Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "001208"
$q1.Range("C2").Value = "诺安低碳经济股票A"
Set-TextValue $q1.Range("D2") "14.51"
Set-TextValue $q1.Range("E2") "82.03"
Set-TextValue $q1.Range("F2") "1.61"
Set-TextValue $q1.Range("G2") "0.2336"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "010349"
$q1.Range("C3").Value = "诺安低碳经济股票C"
Set-TextValue $q1.Range("D3") "3.52"
Set-TextValue $q1.Range("E3") "82.03"
Set-TextValue $q1.Range("F3") "1.61"
Set-TextValue $q1.Range("G3") "0.0567"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "009927"
$q1.Range("C4").Value = "工银瑞信聚利18个月定期开放混合A"
Set-TextValue $q1.Range("D4") "5.54"
Set-TextValue $q1.Range("E4") "23.27"
Set-TextValue $q1.Range("F4") "0.80"
Set-TextValue $q1.Range("G4") "0.0443"
$q1.Range("H4").Value = 9

$q1.Range("A5").Value = 3
Set-TextValue $q1.Range("B5") "009928"
$q1.Range("C5").Value = "工银瑞信聚利18个月定期开放混合C"
Set-TextValue $q1.Range("D5") "0.83"
Set-TextValue $q1.Range("E5") "23.27"
Set-TextValue $q1.Range("F5") "0.80"
Set-TextValue $q1.Range("G5") "0.0066"
$q1.Range("H5").Value = 9

# --- Step 2: append a brand-new "总计" sheet after "2022-Q1" ---------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$headerStyleSource.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$indexStyleSource.Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.34

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.05

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.06

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.06

# Restore the originally-active sheet/tab selection (unrelated to this
# feature, but adding/renaming sheets shifts Excel's "active sheet").
$wb.Worksheets.Item(1).Activate()
